# Updates the cryptos price list (coin name/link swaps + refreshed
# price/volume figures) to match the latest scrape, per the commit
# "Updated cryptos list on Wed Jul 26 20:21:43 UTC 2023 with GitHub Actions".
#
# Column D ("Price") frequently holds values that LOOK numeric
# (e.g. "5.280", "6.000", "0.000007880") but must stay plain text, exactly
# as scraped (Excel would otherwise silently normalise "6.000" -> 6 and
# drop the trailing zeros). To stop Excel's automatic number/date
# detection from mangling those values, column D writes temporarily force
# the cell to Text format ("@") before assigning the value, then restore
# the cell's style to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{ 'B' = 2; 'C' = 3; 'D' = 4; 'E' = 5 }

$edits = @(
    @{ Row = 2; Col = 'D'; Value = '29.461.45' },
    @{ Row = 2; Col = 'E'; Value = '  +0.73%  ' },
    @{ Row = 3; Col = 'D'; Value = '1.880.22' },
    @{ Row = 3; Col = 'E'; Value = '  +1.06%  ' },
    @{ Row = 4; Col = 'D'; Value = '0.9984' },
    @{ Row = 4; Col = 'E'; Value = '  -0.27%  ' },
    @{ Row = 5; Col = 'D'; Value = '0.7139' },
    @{ Row = 5; Col = 'E'; Value = '  +1.99%  ' },
    @{ Row = 6; Col = 'D'; Value = '239.29' },
    @{ Row = 6; Col = 'E'; Value = '  +0.79%  ' },
    @{ Row = 7; Col = 'D'; Value = '0.9994' },
    @{ Row = 7; Col = 'E'; Value = '  -0.21%  ' },
    @{ Row = 8; Col = 'D'; Value = '0.07919' },
    @{ Row = 8; Col = 'E'; Value = '  -2.79%  ' },
    @{ Row = 9; Col = 'D'; Value = '0.3088' },
    @{ Row = 9; Col = 'E'; Value = '  +2.08%  ' },
    @{ Row = 10; Col = 'D'; Value = '25.45' },
    @{ Row = 10; Col = 'E'; Value = '  +9.75%  ' },
    @{ Row = 11; Col = 'D'; Value = '0.08204' },
    @{ Row = 11; Col = 'E'; Value = '  +0.57%  ' },
    @{ Row = 12; Col = 'B'; Value = 'Polkadot' },
    @{ Row = 12; Col = 'C'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Row = 12; Col = 'D'; Value = '5.280' },
    @{ Row = 12; Col = 'E'; Value = '  +2.46%  ' },
    @{ Row = 13; Col = 'B'; Value = 'WrappedEther' },
    @{ Row = 13; Col = 'C'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Row = 13; Col = 'D'; Value = '1.855.75' },
    @{ Row = 13; Col = 'E'; Value = '  +0.07%  ' },
    @{ Row = 14; Col = 'D'; Value = '0.7265' },
    @{ Row = 14; Col = 'E'; Value = '  +2.95%  ' },
    @{ Row = 15; Col = 'D'; Value = '89.59' },
    @{ Row = 15; Col = 'E'; Value = '  +0.65%  ' },
    @{ Row = 16; Col = 'D'; Value = '29.468.97' },
    @{ Row = 16; Col = 'E'; Value = '  +0.71%  ' },
    @{ Row = 17; Col = 'D'; Value = '5.848' },
    @{ Row = 17; Col = 'E'; Value = '  +1.46%  ' },
    @{ Row = 18; Col = 'B'; Value = 'ShibaInu' },
    @{ Row = 18; Col = 'C'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Row = 18; Col = 'D'; Value = '0.000007880' },
    @{ Row = 18; Col = 'E'; Value = '  +0.74%  ' },
    @{ Row = 19; Col = 'B'; Value = 'BitcoinCash' },
    @{ Row = 19; Col = 'C'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Row = 19; Col = 'D'; Value = '242.39' },
    @{ Row = 19; Col = 'E'; Value = '  +2.85%  ' },
    @{ Row = 20; Col = 'D'; Value = '13.39' },
    @{ Row = 20; Col = 'E'; Value = '  +0.58%  ' },
    @{ Row = 21; Col = 'D'; Value = '2.126.67' },
    @{ Row = 21; Col = 'E'; Value = '  +0.87%  ' },
    @{ Row = 22; Col = 'D'; Value = '0.9993' },
    @{ Row = 22; Col = 'E'; Value = '  -0.18%  ' },
    @{ Row = 23; Col = 'D'; Value = '0.9987' },
    @{ Row = 23; Col = 'E'; Value = '  -0.31%  ' },
    @{ Row = 24; Col = 'D'; Value = '7.812' },
    @{ Row = 24; Col = 'E'; Value = '  +5.31%  ' },
    @{ Row = 25; Col = 'B'; Value = 'Monero' },
    @{ Row = 25; Col = 'C'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Row = 25; Col = 'D'; Value = '162.63' },
    @{ Row = 25; Col = 'E'; Value = '  +0.63%  ' },
    @{ Row = 26; Col = 'B'; Value = 'Stellar' },
    @{ Row = 26; Col = 'C'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Row = 26; Col = 'D'; Value = '0.1468' },
    @{ Row = 26; Col = 'E'; Value = '  +2.04%  ' },
    @{ Row = 27; Col = 'D'; Value = '8.995' },
    @{ Row = 27; Col = 'E'; Value = '  +0.53%  ' },
    @{ Row = 28; Col = 'D'; Value = '18.21' },
    @{ Row = 28; Col = 'E'; Value = '  +0.77%  ' },
    @{ Row = 29; Col = 'D'; Value = '1.947' },
    @{ Row = 29; Col = 'E'; Value = '  -0.64%  ' },
    @{ Row = 30; Col = 'D'; Value = '1.357' },
    @{ Row = 30; Col = 'E'; Value = '  -5.23%  ' },
    @{ Row = 31; Col = 'D'; Value = '1.482' },
    @{ Row = 31; Col = 'E'; Value = '  -0.05%  ' },
    @{ Row = 32; Col = 'D'; Value = '4.346' },
    @{ Row = 32; Col = 'E'; Value = '  -1.04%  ' },
    @{ Row = 33; Col = 'D'; Value = '4.105' },
    @{ Row = 33; Col = 'E'; Value = '  +1.32%  ' },
    @{ Row = 34; Col = 'D'; Value = '0.05248' },
    @{ Row = 34; Col = 'E'; Value = '  +1.15%  ' },
    @{ Row = 35; Col = 'D'; Value = '1.197' },
    @{ Row = 35; Col = 'E'; Value = '  +2.59%  ' },
    @{ Row = 36; Col = 'D'; Value = '0.7234' },
    @{ Row = 36; Col = 'E'; Value = '  +2.43%  ' },
    @{ Row = 37; Col = 'D'; Value = '1.003' },
    @{ Row = 37; Col = 'E'; Value = '  +0.23%  ' },
    @{ Row = 38; Col = 'D'; Value = '2.671' },
    @{ Row = 38; Col = 'E'; Value = '  -0.23%  ' },
    @{ Row = 39; Col = 'E'; Value = '  +1.27%  ' },
    @{ Row = 40; Col = 'D'; Value = '2.705' },
    @{ Row = 41; Col = 'D'; Value = '1.177.74' },
    @{ Row = 41; Col = 'E'; Value = '  +3.65%  ' },
    @{ Row = 42; Col = 'D'; Value = '0.9131' },
    @{ Row = 42; Col = 'E'; Value = '  -0.95%  ' },
    @{ Row = 43; Col = 'D'; Value = '6.000' },
    @{ Row = 43; Col = 'E'; Value = '  +1.89%  ' },
    @{ Row = 44; Col = 'B'; Value = 'Aave' },
    @{ Row = 44; Col = 'C'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' },
    @{ Row = 44; Col = 'D'; Value = '72.09' },
    @{ Row = 44; Col = 'E'; Value = '  +2.67%  ' },
    @{ Row = 45; Col = 'B'; Value = 'TheSandbox' },
    @{ Row = 45; Col = 'C'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' },
    @{ Row = 45; Col = 'D'; Value = '0.4327' },
    @{ Row = 45; Col = 'E'; Value = '  +1.65%  ' },
    @{ Row = 46; Col = 'D'; Value = '0.9994' },
    @{ Row = 46; Col = 'E'; Value = '  -0.19%  ' },
    @{ Row = 47; Col = 'D'; Value = '102.55' },
    @{ Row = 47; Col = 'E'; Value = '  +0.34%  ' },
    @{ Row = 48; Col = 'D'; Value = '0.5351' },
    @{ Row = 48; Col = 'E'; Value = '  -1.52%  ' },
    @{ Row = 49; Col = 'D'; Value = '1.777' },
    @{ Row = 49; Col = 'E'; Value = '  +0.76%  ' },
    @{ Row = 50; Col = 'D'; Value = '2.936' },
    @{ Row = 50; Col = 'E'; Value = '  +6.92%  ' },
    @{ Row = 51; Col = 'D'; Value = '9.229' },
    @{ Row = 51; Col = 'E'; Value = '  +0.89%  ' }
)

foreach ($edit in $edits) {
    $cell = $ws.Cells.Item($edit.Row, $colIndex[$edit.Col])
    if ($edit.Col -eq 'D') {
        $cell.NumberFormat = '@'
        $cell.Value = $edit.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $edit.Value
    }
}
